$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F, shifting the existing "District" column (and
# everything after it) one column to the right (F -> G).
$ws.Columns("F:F").Insert()

# New header for the inserted "Address" column.
$ws.Range("F2").Value = "Address"

# Address values for each data row (rows 3-41), derived from the
# school/workplace portion of column B for each teacher.
$addresses = @(
    'Janata Vidyalay ShiraliBhatkal',
    'Assistant TeacherShri KatyayaniHigh School Aversa Ankola',
    'Govt. High School NadumaskeriKumta',
    'Govt. High School GerasoppaHonnavar',
    'G HS BargiKumta',
    'G H S AllankiHonnavar',
    'S P H S MarukeriBhatkal',
    'GunavanteHonnavar',
    'Gibb Girls High SchoolKumta',
    'Assistant TeacherG H S MundalliBhatkal',
    'Janata Vidyalay Kadatoka Honnavar',
    'Assistant TeacherGHS Santeguli',
    'Janata Vidyalaya MudgaKarwar',
    'Assistant TeacherGovt. High SchoolAgragone',
    'H P S JannakadkalHonnavar',
    'Karavali High SchoolBhatakal',
    'KaremmanavarShree Narayan High School NaitesawarWailwad Karwar',
    'S S Comp. PU College KavalakkiHonnavar',
    'Assistant TeacherSharavati High ShoolHonnavar',
    'Hindu High SchoolKarwar',
    'G H S VannalliKumta',
    'Secondary High School AngadiKarwar',
    'Karwar',
    'Assistant TeacherG H S TenginagundiBhatkal',
    'Modern High School HalgaKarwar',
    'Assistant TeacherSatyagrah SmarakVidyalaya ShetgeriAnkola',
    'Shivaji Compsite P U ColleggeChitakulaKarwar',
    'Govt High Scool ChittaraHonnavar',
    'S K P P U College AreangadiHonnavar',
    'Govt. High SchoolTodurKarwar',
    'Govt. High School KalleshwarAnkola',
    'J V Mirjan KodakaniKumta',
    'Sadashivagad High School Sadashivagad',
    'G H S Nelli KeriKumta',
    'G H S BailurBhatkal',
    'Assistant TeacherG H S BandarBhatkal',
    'G H S MankiHonnavar',
    'Janata Vidyalaya',
    'Union High SchoolMajaliKarwar'
)

for ($i = 0; $i -lt $addresses.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 6).Value = $addresses[$i]
}
